$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds the "Förändrad" (changed) date, stored as a date serial.
# Update every data row (2-51) from 45179 (2023-09-10) to 45180 (2023-09-11).
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 3).Value = 45180
}
